$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.921128988265991
$ws.Range("B1").Value = 3.023487329483032
$ws.Range("C1").Value = 3.201297283172607
$ws.Range("D1").Value = 1.058896541595459
$ws.Range("E1").Value = 0.6827820539474487
